$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

# Date header
Replace-Text "2025-12-10 Wednesday" "2025-12-11 Thursday"

# Table cells (order matters to avoid chained replacement collisions:
# 12÷9=1, 3 -> 52÷4=13, 0 must run before 93÷3=31, 0 -> 12÷9=1, 3)
Replace-Text "74÷8=9, 2" "15÷9=1, 6"
Replace-Text "67÷8=8, 3" "11÷2=5, 1"
Replace-Text "22÷9=2, 4" "88÷7=12, 4"
Replace-Text "86÷9=9, 5" "30÷4=7, 2"
Replace-Text "74÷2=37, 0" "35÷4=8, 3"
Replace-Text "54÷2=27, 0" "68÷4=17, 0"
Replace-Text "83÷9=9, 2" "91÷9=10, 1"
Replace-Text "17÷9=1, 8" "10÷7=1, 3"
Replace-Text "98÷2=49, 0" "87÷9=9, 6"
Replace-Text "50÷6=8, 2" "90÷4=22, 2"
Replace-Text "12÷9=1, 3" "52÷4=13, 0"
Replace-Text "10÷9=1, 1" "59÷5=11, 4"
Replace-Text "36÷9=4, 0" "15÷9=1, 6"
Replace-Text "29÷2=14, 1" "65÷5=13, 0"
Replace-Text "28÷3=9, 1" "79÷4=19, 3"
Replace-Text "63÷2=31, 1" "11÷3=3, 2"
Replace-Text "58÷6=9, 4" "42÷2=21, 0"
Replace-Text "29÷8=3, 5" "13÷2=6, 1"
Replace-Text "92÷2=46, 0" "31÷6=5, 1"
Replace-Text "93÷3=31, 0" "12÷9=1, 3"
Replace-Text "45÷9=5, 0" "42÷5=8, 2"
Replace-Text "34÷7=4, 6" "25÷9=2, 7"
Replace-Text "15÷8=1, 7" "20÷6=3, 2"
Replace-Text "38÷4=9, 2" "44÷4=11, 0"
Replace-Text "13÷4=3, 1" "62÷6=10, 2"
